# Applies the "added number density capability; fixed minor bugs" edit:
#  - sheet "exp": insert 3 new experiment rows (lowO2 / midO2 / highO2)
#    before the moshammer row, keeping outcome/plot/plot columns.
#  - sheet "mech": rename the v37/v38/v39 mechanism rows to v51/v52/v53,
#    relabel x10 -> /3 row to x10, x10 row to x30, and delete the old
#    v40 (/10) row entirely (table shrinks from 6 rows to 5 rows).

$wb = $excel.ActiveWorkbook

$expSheet = $wb.Worksheets.Item("exp")
$mechSheet = $wb.Worksheets.Item("mech")

# --- sheet "exp": insert 3 new rows above the moshammer row (row 8) ---
$moshammerRow = 8
$expSheet.Rows.Item($moshammerRow).Insert()
$expSheet.Rows.Item($moshammerRow).Insert()
$expSheet.Rows.Item($moshammerRow).Insert()

$expSheet.Cells.Item($moshammerRow, 1).Value = "couch_2022_dme_lowO2.xlsx"
$expSheet.Cells.Item($moshammerRow, 2).Value = "outcome"
$expSheet.Cells.Item($moshammerRow, 3).Value = "plot"
$expSheet.Cells.Item($moshammerRow, 4).Value = "plot"

$expSheet.Cells.Item($moshammerRow + 1, 1).Value = "couch_2022_dme_midO2.xlsx"
$expSheet.Cells.Item($moshammerRow + 1, 2).Value = "outcome"
$expSheet.Cells.Item($moshammerRow + 1, 3).Value = "plot"
$expSheet.Cells.Item($moshammerRow + 1, 4).Value = "plot"

$expSheet.Cells.Item($moshammerRow + 2, 1).Value = "couch_2022_dme_highO2.xlsx"
$expSheet.Cells.Item($moshammerRow + 2, 2).Value = "outcome"
$expSheet.Cells.Item($moshammerRow + 2, 3).Value = "plot"
$expSheet.Cells.Item($moshammerRow + 2, 4).Value = "plot"

[void]$expSheet.Range("C15").Select()
$expSheet.Application.ActiveWindow.DisplayFormulas = $true

# --- sheet "mech": update mechanism rows ---
$mechSheet.Cells.Item(3, 1).Value = "dme_couch_v51.cti"
$mechSheet.Cells.Item(3, 3).Value = "x3"

$mechSheet.Cells.Item(4, 1).Value = "dme_couch_v52.cti"
$mechSheet.Cells.Item(4, 3).Value = "x10"

$mechSheet.Cells.Item(5, 1).Value = "dme_couch_v53.cti"
$mechSheet.Cells.Item(5, 3).Value = "x30"

# remove the old 6th row (dme_couch_v40.cti / /10)
$mechSheet.Rows.Item(6).Delete()

[void]$mechSheet.Range("A6").Select()
